# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rows 16-31 hold the "Periodo Mora" (col E), "Valor Mora" (col F) and
# "Salario Basico" (col G) table. The period list is refreshed to the new
# ascending sequence (2312, 2401 .. 2503) and the corresponding amounts are
# updated to match.

$wb = $excel.ActiveWorkbook

if ($wb.Worksheets | Where-Object { $_.Name -eq "Hoja1" }) {
    $ws = $wb.Worksheets.Item("Hoja1")
} else {
    $ws = $wb.ActiveSheet
}

# period, Valor Mora (F), Salario Basico (G) for rows 16..31
$rows = @(
    @{ Row = 16; Periodo = "2312"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 17; Periodo = "2401"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 18; Periodo = "2402"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 19; Periodo = "2403"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 20; Periodo = "2404"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 21; Periodo = "2405"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 22; Periodo = "2406"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 23; Periodo = "2407"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 24; Periodo = "2408"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 25; Periodo = "2409"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 26; Periodo = "2410"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 27; Periodo = "2411"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 28; Periodo = "2412"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 29; Periodo = "2501"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 30; Periodo = "2502"; ValorMora = 46400;  Salario = 1160000 },
    @{ Row = 31; Periodo = "2503"; ValorMora = 27840;  Salario = 1160000 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("E$r").Value = $item.Periodo
    $ws.Range("F$r").Value = $item.ValorMora
    $ws.Range("G$r").Value = $item.Salario
}
